$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Date-text cells that must remain literal text (not auto-converted to a
# date serial number) even though they look like dates.
$dateUpdates = @{
    "B2" = "2024.02.13"
    "B3" = "2024.02.17"
    "B4" = "2024.02.19"
    "B5" = "2024.03.02"
    "B6" = "2024.03.16"
    "B7" = "2024.03.23"
    "B8" = "2024.04.04"
    "B9" = "2024.05.18"
}

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    foreach ($addr in $dateUpdates.Keys) {
        $cell = $ws.Range($addr)
        # Force text interpretation so the dotted date string isn't
        # reinterpreted as a date serial, then restore the default
        # (unstyled) cell style so no stray formatting is introduced.
        $cell.NumberFormat = "@"
        $cell.Value = $dateUpdates[$addr]
        $cell.Style = "Normal"
    }

    $ws.Range("F3").Value = 1868

    $ws.Range("F4").Value = 353

    $ws.Range("F5").Value = 1139
    $ws.Range("G5").Value = 49.5

    $ws.Range("F6").Value = 1136

    $ws.Range("F8").Value = 5964
}
